$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 6 (the stray formatted-but-empty row)
$ws.Rows("6:6").Delete()

# --- Phase 1: set NEW header text in the exact order needed for sharedStrings append order ---
$ws.Range("O1").Value = "ECORP_COUNTY"
$ws.Range("N1").Value = "ECORP_STATE"
$ws.Range("M1").Value = "ECORP_BUSINESS_TYPE"
$ws.Range("L1").Value = "ECORP_FORMATION_DATE"
$ws.Range("K1").Value = "ECORP_STATUS"
$ws.Range("J1").Value = "ECORP_ENTITY_TYPE"
$ws.Range("I1").Value = "ECORP_ENTITY_ID_S"
$ws.Range("H1").Value = "ECORP_NAME_S"
$ws.Range("G1").Value = "ECORP_TYPE"
$ws.Range("F1").Value = "ECORP_SEARCH_NAME"
$ws.Range("P1").Value = "ECORP_COMMENTS"

# --- Phase 2: set remaining headers (reuse of already-existing shared strings; order does not matter) ---
$ws.Range("A1").Value = "FULL_ADDRESS"
$ws.Range("B1").Value = "COUNTY"
$ws.Range("C1").Value = "Owner_Ownership"
$ws.Range("D1").Value = "ECORP_INDEX_#"
$ws.Range("E1").Value = "OWNER_TYPE"
$ws.Range("Q1").Value = "StatutoryAgent1_Name"
$ws.Range("R1").Value = "StatutoryAgent1_Address"
$ws.Range("S1").Value = "StatutoryAgent1_Phone"
$ws.Range("T1").Value = "StatutoryAgent1_Mail"
$ws.Range("U1").Value = "StatutoryAgent2_Name"
$ws.Range("V1").Value = "StatutoryAgent2_Address"
$ws.Range("W1").Value = "StatutoryAgent2_Phone"
$ws.Range("X1").Value = "StatutoryAgent2_Mail"
$ws.Range("Y1").Value = "StatutoryAgent3_Name"
$ws.Range("Z1").Value = "StatutoryAgent3_Address"
$ws.Range("AA1").Value = "StatutoryAgent3_Phone"
$ws.Range("AB1").Value = "StatutoryAgent3_Mail"
$ws.Range("AC1").Value = "Manager1_Name"
$ws.Range("AD1").Value = "Manager1_Address"
$ws.Range("AE1").Value = "Manager1_Phone"
$ws.Range("AF1").Value = "Manager1_Mail"
$ws.Range("AG1").Value = "Manager2_Name"
$ws.Range("AH1").Value = "Manager2_Address"
$ws.Range("AI1").Value = "Manager2_Phone"
$ws.Range("AJ1").Value = "Manager2_Mail"
$ws.Range("AK1").Value = "Manager3_Name"
$ws.Range("AL1").Value = "Manager3_Address"
$ws.Range("AM1").Value = "Manager3_Phone"
$ws.Range("AN1").Value = "Manager3_Mail"
$ws.Range("AO1").Value = "Manager4_Name"
$ws.Range("AP1").Value = "Manager4_Address"
$ws.Range("AQ1").Value = "Manager4_Phone"
$ws.Range("AR1").Value = "Manager4_Mail"
$ws.Range("AS1").Value = "Manager5_Name"
$ws.Range("AT1").Value = "Manager5_Address"
$ws.Range("AU1").Value = "Manager5_Phone"
$ws.Range("AV1").Value = "Manager5_Mail"
$ws.Range("AW1").Value = "Manager/Member1_Name"
$ws.Range("AX1").Value = "Manager/Member1_Address"
$ws.Range("AY1").Value = "Manager/Member1_Phone"
$ws.Range("AZ1").Value = "Manager/Member1_Mail"
$ws.Range("BA1").Value = "Manager/Member2_Name"
$ws.Range("BB1").Value = "Manager/Member2_Address"
$ws.Range("BC1").Value = "Manager/Member2_Phone"
$ws.Range("BD1").Value = "Manager/Member2_Mail"
$ws.Range("BE1").Value = "Manager/Member3_Name"
$ws.Range("BF1").Value = "Manager/Member3_Address"
$ws.Range("BG1").Value = "Manager/Member3_Phone"
$ws.Range("BH1").Value = "Manager/Member3_Mail"
$ws.Range("BI1").Value = "Manager/Member4_Name"
$ws.Range("BJ1").Value = "Manager/Member4_Address"
$ws.Range("BK1").Value = "Manager/Member4_Phone"
$ws.Range("BL1").Value = "Manager/Member4_Mail"
$ws.Range("BM1").Value = "Manager/Member5_Name"
$ws.Range("BN1").Value = "Manager/Member5_Address"
$ws.Range("BO1").Value = "Manager/Member5_Phone"
$ws.Range("BP1").Value = "Manager/Member5_Mail"
$ws.Range("BQ1").Value = "Member1_Name"
$ws.Range("BR1").Value = "Member1_Address"
$ws.Range("BS1").Value = "Member1_Phone"
$ws.Range("BT1").Value = "Member1_Mail"
$ws.Range("BU1").Value = "Member2_Name"
$ws.Range("BV1").Value = "Member2_Address"
$ws.Range("BW1").Value = "Member2_Phone"
$ws.Range("BX1").Value = "Member2_Mail"
$ws.Range("BY1").Value = "Member3_Name"
$ws.Range("BZ1").Value = "Member3_Address"
$ws.Range("CA1").Value = "Member3_Phone"
$ws.Range("CB1").Value = "Member3_Mail"
$ws.Range("CC1").Value = "Member4_Name"
$ws.Range("CD1").Value = "Member4_Address"
$ws.Range("CE1").Value = "Member4_Phone"
$ws.Range("CF1").Value = "Member4_Mail"
$ws.Range("CG1").Value = "Member5_Name"
$ws.Range("CH1").Value = "Member5_Address"
$ws.Range("CI1").Value = "Member5_Phone"
$ws.Range("CJ1").Value = "Member5_Mail"
$ws.Range("CK1").Value = "IndividualName1"
$ws.Range("CL1").Value = "IndividualName2"
$ws.Range("CM1").Value = "IndividualName3"
$ws.Range("CN1").Value = "IndividualName4"
$ws.Range("CO1").Value = "ECORP_URL"

# --- Apply cell formatting ---
# Q1:CN1 -> bold, NO fill, full thin border (new style introduced by the edit)
$noFillRange = $ws.Range("Q1:CN1")
$noFillRange.Font.Bold = $true
$noFillRange.Font.Size = 11
$noFillRange.Borders.LineStyle = 1
$noFillRange.Borders.Weight = 2
$noFillRange.HorizontalAlignment = -4108
$noFillRange.VerticalAlignment = -4160

# Reset selection / view state to match target
$ws.Range("B3").Select()
